# Trading update: 2026-02-17 07:58:56
# Appends a new trade row (row 15 / Trade # 14) to both the "All Trades"
# and "MarketMaking" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 15

    # A: Trade #
    $ws.Cells.Item($row, 1).Value = 14

    # B: Date -- force text so the date-like string is not reinterpreted
    # as a date serial number, then drop the quote-prefix style again so
    # no extra cell style is left behind.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).Style = "Normal"

    # C: Time
    $ws.Cells.Item($row, 3).Value = "07:58:50"

    # D: Strategy
    $ws.Cells.Item($row, 4).Value = "MarketMaking"

    # E: Side
    $ws.Cells.Item($row, 5).Value = "DOWN"

    # F: Entry Price
    $ws.Cells.Item($row, 6).Value = 0.93

    # G: Exit Price -- empty (trade still open). Use a lone quote prefix to
    # materialize an empty text cell, then clear the style it introduces.
    $ws.Cells.Item($row, 7).Value = "'"
    $ws.Cells.Item($row, 7).Style = "Normal"

    # H: Status
    $ws.Cells.Item($row, 8).Value = "OPEN"

    # I: P&L %
    $ws.Cells.Item($row, 9).Value = 0

    # J: P&L $
    $ws.Cells.Item($row, 10).Value = 0

    # K: Capital After
    $ws.Cells.Item($row, 11).Value = 99.97999999999999

    # L: Entry Slippage (bps)
    $ws.Cells.Item($row, 12).Value = 0

    # M: Exit Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0

    # N: Confidence
    $ws.Cells.Item($row, 14).Value = 0.6

    # O: Entry Reason
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"

    # P: Exit Reason -- empty (trade still open)
    $ws.Cells.Item($row, 16).Value = "'"
    $ws.Cells.Item($row, 16).Style = "Normal"

    # Q: Duration (min)
    $ws.Cells.Item($row, 17).Value = 0
}
